$wb = $excel.ActiveWorkbook
$list1 = $wb.Worksheets.Item("List1")
$ws = $wb.Worksheets.Add($null, $list1)

$ws.Range("C5").Value = "Day 1 "
$ws.Range("D5").Value = "Day 2 "
$ws.Range("E5").Value = "Day 3 "
$ws.Range("F5").Value = "Day 4 "

$ws.Range("B6").Value = "Actual streamline"
$ws.Range("C6").Value = 140
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 45
$ws.Range("F6").Value = 5

$ws.Range("B7").Value = "Desired streamline"
$ws.Range("C7").Value = 140
$ws.Range("D7").Value = 95
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 0

$chartObj = $ws.ChartObjects().Add(100, 100, 400, 300)
$chart = $chartObj.Chart
$chart.ChartType = 4

$chart.SeriesCollection().NewSeries()
$s1 = $chart.SeriesCollection(1)
$s1.Name = "=Sheet1!`$B`$6"
$s1.Values = "=Sheet1!`$C`$6:`$F`$6"
$s1.XValues = "=Sheet1!`$C`$5:`$F`$5"

$chart.SeriesCollection().NewSeries()
$s2 = $chart.SeriesCollection(2)
$s2.Name = "=Sheet1!`$B`$7"
$s2.Values = "=Sheet1!`$C`$7:`$F`$7"
$s2.XValues = "=Sheet1!`$C`$5:`$F`$5"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sprint 2"
$chart.HasLegend = $true
$chart.Legend.Position = -4107

Write-Host "chart added"
